# fix: createDataSet 경로 수정
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("F2").Value = "EGC3035"
$ws.Range("H2").Value = "바이오의료산업경영"
$ws.Range("I2").Value = "권지연"
$ws.Range("K2").Value = "B0"
$ws.Range("T2").Value = "Introduction to Health Industry Management"

# Row 3 updates
$ws.Range("F3").Value = "EGC4023"
$ws.Range("H3").Value = "인간과우주"
$ws.Range("I3").Value = "이관수"
$ws.Range("K3").Value = "A0"
$ws.Range("Q3").Value = ""
$ws.Range("T3").Value = "Humanity and Universe"
